# "Ajuste con factor de correccion" - update the aforo (capacity) figure
# for "Plaza de toros" (B3) with the corrected value, and leave the
# selection on the edited cell, matching the authored change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 7000
$ws.Range("B3").Select()
